$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 12 at the bottom of the results table. Inserting (rather
# than just writing into a blank row) makes Excel carry the formatting of
# the row above down into the new row, same as a user pressing Tab/Enter
# after the last row of the table.
$ws.Rows(12).Insert()

# Row 11's "D" column (SpotId) happens to use the default style, but the
# new row's SpotId cell should pick up the numeric-ish style (s=2) used by
# the other data columns - copy that formatting over explicitly.
$ws.Range("C11").Copy()
$ws.Range("D12").PasteSpecial(-4122)

# New fake result row: firstNameNumber, secondNameNumber, PerformanceId,
# SpotId, winnerId, comments1, comments2
$ws.Range("A12").Value = "Hastings.286"
$ws.Range("B12").Value = "Mungo.80"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = "I2"
$ws.Range("E12").Value = "Hastings.286"
$ws.Range("F12").Value = "Great 8 to 5. Triple box TTR could use some work. Needs more intensity. Sloopy unstable. Playing inconsistent at best"
$ws.Range("G12").Value = "Poor marching and playing"

# Match the row height used by the rest of the data rows.
$ws.Rows(12).RowHeight = 18

# Leave the selection where the user would land after finishing the row.
$null = $ws.Range("C13").Select()
